# Fruta / hortaliza, semanal
# Insert 3 new weekly rows for "Macroferia Regional de Talca - Frutilla"
# right before the existing row 260, shifting the rest of the block down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at position 260 (pushes old rows 260:284 -> 263:287)
$ws.Range("A260:A262").EntireRow.Insert()

# --- New row 260 ---
$ws.Cells.Item(260, 1).Value = 5
$ws.Cells.Item(260, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(260, 3).Value = "Maule"
$ws.Cells.Item(260, 4).Value = 44461
$ws.Cells.Item(260, 5).Value = 7
$ws.Cells.Item(260, 6).Value = "Fruta"
$ws.Cells.Item(260, 7).Value = 100101
$ws.Cells.Item(260, 8).Value = "Berries"
$ws.Cells.Item(260, 9).Value = 100112025
$ws.Cells.Item(260, 10).Value = "Frutilla"
$ws.Cells.Item(260, 11).Value = "Sin especificar"
$ws.Cells.Item(260, 12).Value = "Especial"
$ws.Cells.Item(260, 13).Value = 50
$ws.Cells.Item(260, 14).Value = 20000
$ws.Cells.Item(260, 15).Value = 20000
$ws.Cells.Item(260, 16).Value = 20000
$ws.Cells.Item(260, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(260, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(260, 19).Value = 2857
$ws.Cells.Item(260, 20).Value = 7

# --- New row 261 ---
$ws.Cells.Item(261, 1).Value = 5
$ws.Cells.Item(261, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(261, 3).Value = "Maule"
$ws.Cells.Item(261, 4).Value = 44461
$ws.Cells.Item(261, 5).Value = 7
$ws.Cells.Item(261, 6).Value = "Fruta"
$ws.Cells.Item(261, 7).Value = 100101
$ws.Cells.Item(261, 8).Value = "Berries"
$ws.Cells.Item(261, 9).Value = 100112025
$ws.Cells.Item(261, 10).Value = "Frutilla"
$ws.Cells.Item(261, 11).Value = "Sin especificar"
$ws.Cells.Item(261, 12).Value = "Primera"
$ws.Cells.Item(261, 13).Value = 40
$ws.Cells.Item(261, 14).Value = 17000
$ws.Cells.Item(261, 15).Value = 17000
$ws.Cells.Item(261, 16).Value = 17000
$ws.Cells.Item(261, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(261, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(261, 19).Value = 2429
$ws.Cells.Item(261, 20).Value = 7

# --- New row 262 ---
$ws.Cells.Item(262, 1).Value = 5
$ws.Cells.Item(262, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(262, 3).Value = "Maule"
$ws.Cells.Item(262, 4).Value = 44461
$ws.Cells.Item(262, 5).Value = 7
$ws.Cells.Item(262, 6).Value = "Fruta"
$ws.Cells.Item(262, 7).Value = 100101
$ws.Cells.Item(262, 8).Value = "Berries"
$ws.Cells.Item(262, 9).Value = 100112025
$ws.Cells.Item(262, 10).Value = "Frutilla"
$ws.Cells.Item(262, 11).Value = "Sin especificar"
$ws.Cells.Item(262, 12).Value = "Segunda"
$ws.Cells.Item(262, 13).Value = 30
$ws.Cells.Item(262, 14).Value = 12000
$ws.Cells.Item(262, 15).Value = 12000
$ws.Cells.Item(262, 16).Value = 12000
$ws.Cells.Item(262, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(262, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(262, 19).Value = 1714
$ws.Cells.Item(262, 20).Value = 7
